$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address paragraph "919 Story Road, San Jose CA 95122"
#    into two separate paragraphs: "919 Story Road" and "San Jose, CA 95122".
#    (Only the first occurrence -- the mailing address near the top of the
#    letter -- changes; the later "PROPERTY ADDRESS" occurrence stays intact.)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "919 Story Road, San Jose CA 95122`r") {
        $p.Range.Text = "919 Story Road`r"
        $newPara = $p.Next()
        $newPara.Range.Text = "San Jose, CA 95122"
        break
    }
}

# 3. Remove the empty "NoSpacing" paragraph directly after "Board of Directors"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors`r") {
        $next = $p.Next()
        $next.Range.Delete()
        break
    }
}
